# Updated cryptos list on Sat Mar 16 07:41:54 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be treated as text while we write values that
# look numeric (e.g. "611.73"), matching the original inlineStr cells, then
# restore the default "Normal" style so no stray number-format survives on
# the cells.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "69.326.83"
$ws.Range("E2").Value = "  +1.23%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.728.59"
$ws.Range("E3").Value = "  -0.62%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.12%  "

# Row 5 - BNB
$ws.Range("D5").Value = "611.73"
$ws.Range("E5").Value = "  +4.94%  "

# Row 6 - Solana
$ws.Range("D6").Value = "192.24"
$ws.Range("E6").Value = "  +9.08%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.638"
$ws.Range("E7").Value = "  +0.26%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.02%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.730"
$ws.Range("E9").Value = "  +1.25%  "

# Row 10 - now Dogecoin (was Avalanche)
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.162"
$ws.Range("E10").Value = "  -3.29%  "

# Row 11 - now Avalanche (was Dogecoin)
$ws.Range("B11").Value = "Avalanche"
$ws.Range("C11").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D11").Value = "60.60"
$ws.Range("E11").Value = "  +12.42%  "

# Row 12 - ShibaInu
$ws.Range("D12").Value = "0.0000291"
$ws.Range("E12").Value = "  -3.69%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "10.69"
$ws.Range("E13").Value = "  -1.42%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "4.316.91"
$ws.Range("E14").Value = "  -0.31%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "3.727.56"
$ws.Range("E15").Value = "  -0.69%  "

# Row 16 - Chainlink
$ws.Range("D16").Value = "19.49"
$ws.Range("E16").Value = "  -0.52%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  +0.20%  "

# Row 18 - TRON
$ws.Range("E18").Value = "  -0.21%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "12.97"
$ws.Range("E19").Value = "  -1.61%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "69.156.58"
$ws.Range("E20").Value = "  +1.22%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "413.22"
$ws.Range("E21").Value = "  -0.31%  "

# Row 22 - PancakeSwap
$ws.Range("D22").Value = "4.58"
$ws.Range("E22").Value = "  +0.20%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "89.86"
$ws.Range("E23").Value = "  +0.41%  "

# Row 24 - ImmutableX
$ws.Range("E24").Value = "  -1.54%  "

# Row 25 - InternetComputer(DFINITY)
$ws.Range("D25").Value = "12.90"
$ws.Range("E25").Value = "  -0.94%  "

# Row 26 - RenderToken
$ws.Range("D26").Value = "10.92"
$ws.Range("E26").Value = "  +0.31%  "

# Row 27 - Toncoin
$ws.Range("E27").Value = "  -1.41%  "

# Row 28 - LEO
$ws.Range("D28").Value = "6.06"
$ws.Range("E28").Value = "  +1.26%  "

# Row 29 - Filecoin
$ws.Range("D29").Value = "9.72"
$ws.Range("E29").Value = "  +0.19%  "

# Row 30 - EthereumClassic
$ws.Range("D30").Value = "33.17"
$ws.Range("E30").Value = "  -0.60%  "

# Row 31 - NEARProtocol
$ws.Range("D31").Value = "7.79"
$ws.Range("E31").Value = "  -3.29%  "

# Row 32 - Cosmos
$ws.Range("D32").Value = "12.86"
$ws.Range("E32").Value = "  -0.05%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "0.125"
$ws.Range("E33").Value = "  +4.57%  "

# Row 34 - InjectiveProtocol
$ws.Range("D34").Value = "45.94"
$ws.Range("E34").Value = "  +3.76%  "

# Row 35 - Bittensor
$ws.Range("D35").Value = "637.03"
$ws.Range("E35").Value = "  +3.38%  "

# Row 36 - OKB
$ws.Range("D36").Value = "66.07"
$ws.Range("E36").Value = "  +0.28%  "

# Row 37 - TheGraph
$ws.Range("D37").Value = "0.419"
$ws.Range("E37").Value = "  +3.07%  "

# Row 38 - PEPE
$ws.Range("D38").Value = "0.0₃0831"
$ws.Range("E38").Value = "  -11.12%  "

# Row 39 - Dai
$ws.Range("E39").Value = "  -0.17%  "

# Row 40 - FirstDigitalUSD
$ws.Range("E40").Value = "  +0.22%  "

# Row 41 - Kaspa
$ws.Range("E41").Value = "  +3.06%  "

# Row 42 - ThetaToken
$ws.Range("D42").Value = "3.06"
$ws.Range("E42").Value = "  -1.39%  "

# Row 43 - VeChain
$ws.Range("D43").Value = "0.0450"
$ws.Range("E43").Value = "  +0.14%  "

# Row 44 - Fetch.AI
$ws.Range("D44").Value = "2.64"
$ws.Range("E44").Value = "  -0.22%  "

# Row 45 - Stellar
$ws.Range("E45").Value = "  +2.59%  "

# Row 46 - Maker
$ws.Range("D46").Value = "2.890.34"
$ws.Range("E46").Value = "  +4.77%  "

# Row 47 - THORChain
$ws.Range("E47").Value = "  -3.44%  "

# Row 48 - WEMIXToken
$ws.Range("E48").Value = "  +0.34%  "

# Row 49 - Monero
$ws.Range("D49").Value = "144.00"
$ws.Range("E49").Value = "  +0.27%  "

# Row 50 - ApeXProtocol
$ws.Range("D50").Value = "3.10"
$ws.Range("E50").Value = "  -0.90%  "

# Row 51 - Stacks
$ws.Range("D51").Value = "2.79"
$ws.Range("E51").Value = "  -0.27%  "

# Restore default styling on the Price column so no stray number-format
# sticks around on cells we touched.
$ws.Range("D2:D51").Style = "Normal"
